# Update the "F" column (view/visit counters) figures on each sheet to
# match the newly generated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 3276
$ws.Range("F4").Value  = 2015
$ws.Range("F5").Value  = 273
$ws.Range("F6").Value  = 116
$ws.Range("F7").Value  = 3100
$ws.Range("F8").Value  = 620
$ws.Range("F10").Value = 43
$ws.Range("F13").Value = 157
$ws.Range("F14").Value = 157
$ws.Range("F15").Value = 10239
$ws.Range("F17").Value = 239
$ws.Range("F20").Value = 8128
$ws.Range("F21").Value = 12738
$ws.Range("F24").Value = 32
$ws.Range("F26").Value = 402
$ws.Range("F28").Value = 9
$ws.Range("F29").Value = 428
$ws.Range("F30").Value = 2832
$ws.Range("F31").Value = 261
$ws.Range("F33").Value = 8026
$ws.Range("F34").Value = 1640
$ws.Range("F38").Value = 4631
$ws.Range("F39").Value = 1483
$ws.Range("F43").Value = 645
$ws.Range("F44").Value = 4

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value  = 1204
$ws.Range("F13").Value = 73
$ws.Range("F15").Value = 16

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 24

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 3276
$ws.Range("F6").Value  = 2015
$ws.Range("F8").Value  = 273
$ws.Range("F9").Value  = 24
$ws.Range("F10").Value = 3100
$ws.Range("F12").Value = 620
$ws.Range("F13").Value = 43
$ws.Range("F16").Value = 157
$ws.Range("F17").Value = 157
$ws.Range("F18").Value = 10239
$ws.Range("F19").Value = 239
$ws.Range("F22").Value = 8128
$ws.Range("F23").Value = 12738
$ws.Range("F25").Value = 32
$ws.Range("F30").Value = 9
$ws.Range("F31").Value = 2832
$ws.Range("F34").Value = 261
$ws.Range("F36").Value = 8026
$ws.Range("F40").Value = 4631
$ws.Range("F47").Value = 645
